$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws4 = $wb.Worksheets.Item(4)

# ---------------------------------------------------------------------------
# 1. Populate cell values. The order below controls the order in which new
#    strings are appended to the shared-strings table, so it must mirror the
#    order the strings first appear in the target workbook.
# ---------------------------------------------------------------------------
$ws1.Range("K1").Value2 = "Skrót"
$ws1.Range("L4").Value2 = "Sum Square Difference"
$ws1.Range("L5").Value2 = "Cross - Correlation"
$ws1.Range("L2").Value2 = "Mean Square Error"
$ws1.Range("M1").Value2 = "Opis"
$ws1.Range("L1").Value2 = "Rozwinięcie"
$ws1.Range("M2").Value2 = "Im mniejsza wartość MSE, tym obrazy są bardziej podobne."
$ws1.Range("L3").Value2 = "Mean Absolute Error"
$ws1.Range("M3").Value2 = "Ocenia, jaka jest przeciętna różnica absolutna między odpowiadającymi sobie komórkami danych."
$ws1.Range("M4").Value2 = "Im mniejsza wartość SSD, tym większe podobieństwo obrazów. Wrażliwy na duże różnice w wartościach (podnoszone do kwadratu)."
$ws1.Range("M5").Value2 = "Miara podobieństwa między obrazami, uwzględniająca przesunięcie (lub dopasowanie) jednego względem drugiego. Im większa wartość tym bardziej dopasowane."

$ws1.Range("K2").Value2 = "MSE"
$ws1.Range("K3").Value2 = "MAE"
$ws1.Range("K4").Value2 = "SSD"
$ws1.Range("K5").Value2 = "CC"

# ---------------------------------------------------------------------------
# 2. Apply formatting. Order matters here too: each distinct combination of
#    font/border creates a new cell style the first time it is used, so we
#    touch the cells in the same order the styles appear in the target file.
# ---------------------------------------------------------------------------

# Header row K1:M1 -> bold font + thin border (style used for "Skrót" / "Rozwinięcie" / "Opis")
$ws1.Range("K1").Borders.LineStyle = 1
$ws1.Range("K1").Font.Bold = $true
$ws1.Range("L1").Borders.LineStyle = 1
$ws1.Range("L1").Font.Bold = $true
$ws1.Range("M1").Borders.LineStyle = 1
$ws1.Range("M1").Font.Bold = $true

# Plain bordered cells (default font) for the bulk of the table
$ws1.Range("K2").Borders.LineStyle = 1
$ws1.Range("L2").Borders.LineStyle = 1
$ws1.Range("K3").Borders.LineStyle = 1
$ws1.Range("L3").Borders.LineStyle = 1
$ws1.Range("K4").Borders.LineStyle = 1
$ws1.Range("L4").Borders.LineStyle = 1
$ws1.Range("K5").Borders.LineStyle = 1
$ws1.Range("L5").Borders.LineStyle = 1
$ws1.Range("M3").Borders.LineStyle = 1
$ws1.Range("M4").Borders.LineStyle = 1
$ws1.Range("M5").Borders.LineStyle = 1

# M2 -> explicit black font + thin border
$ws1.Range("M2").Borders.LineStyle = 1
$ws1.Range("M2").Font.Color = 0

# ---------------------------------------------------------------------------
# 3. Column widths for the newly-populated columns (best-fit-like sizing).
# ---------------------------------------------------------------------------
$ws1.Columns.Item(2).EntireColumn.AutoFit()
$ws1.Columns.Item(3).EntireColumn.AutoFit()
$ws1.Columns.Item(4).EntireColumn.AutoFit()
$ws1.Columns.Item(5).EntireColumn.AutoFit()
$ws1.Columns.Item(7).EntireColumn.AutoFit()
$ws1.Columns.Item(12).EntireColumn.AutoFit()
$ws1.Columns.Item(13).EntireColumn.AutoFit()

# ---------------------------------------------------------------------------
# 4. Switch the active sheet/selection from METIS_KOMORA_1m back to
#    METIS_AINFO_1m, matching the updated workbook view state.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("K21").Select()
